# Courchevel - Soutenance finale : mise à jour du schéma d'accréditation
# et nettoyage des captures d'écran MCD (fond blanc -> transparent).

$p = $ppt.ActivePresentation
$EMU = 12700  # EMU per point

# ---------------------------------------------------------------------
# Slide 14 : repositionnement de quelques ellipses / connecteurs
# ---------------------------------------------------------------------
$s14 = $p.Slides.Item(14)

$shp = $s14.Shapes.Item("Ellipse 11")
$shp.Left = 3538358 / $EMU
$shp.Top  = 3861048 / $EMU

$shp = $s14.Shapes.Item("Connecteur droit 14")
$shp.Width  = 1572735 / $EMU
$shp.Height = 724611 / $EMU

$shp = $s14.Shapes.Item("Ellipse 17")
$shp.Left = 3491880 / $EMU
$shp.Top  = 5169390 / $EMU

$shp = $s14.Shapes.Item("Connecteur droit 27")
$shp.Width  = 1526257 / $EMU
$shp.Height = 1951887 / $EMU

# ---------------------------------------------------------------------
# Slide 15 : repositionnement de quelques ellipses / connecteurs
# et texte de l'ellipse "Ellipse 52"
# ---------------------------------------------------------------------
$s15 = $p.Slides.Item(15)

$shp = $s15.Shapes.Item("Ellipse 9")
$shp.Left = 3779912 / $EMU
$shp.Top  = 1960180 / $EMU

$shp = $s15.Shapes.Item("Ellipse 10")
$shp.Left = 4355976 / $EMU
$shp.Top  = 3239438 / $EMU

$shp = $s15.Shapes.Item("Ellipse 11")
$shp.Left = 3779912 / $EMU
$shp.Top  = 4437112 / $EMU

$shp = $s15.Shapes.Item("Connecteur droit 12")
$shp.Left   = 1349608 / $EMU
$shp.Top    = 2530235 / $EMU
$shp.Width  = 2430304 / $EMU
$shp.Height = 1086648 / $EMU
$shp.VerticalFlip = -1

$shp = $s15.Shapes.Item("Connecteur droit 13")
$shp.Width  = 3006368 / $EMU
$shp.Height = 162615 / $EMU

$shp = $s15.Shapes.Item("Connecteur droit 14")
$shp.Left   = 1349608 / $EMU
$shp.Top    = 3616883 / $EMU
$shp.Width  = 2430304 / $EMU
$shp.Height = 1405351 / $EMU
$shp.VerticalFlip = 0

$shp = $s15.Shapes.Item("Ellipse 52")
$tr = $shp.TextFrame.TextRange
# Remplace le 2e run (qui porte déjà smtClean="0") par le nouveau texte,
# puis vide le 1er run pour ne garder qu'un seul <a:r>.
$tail = $tr.Characters(19, 14)
$tail.Text = "Imprimer une accréditation"
$head = $tr.Characters(1, 18)
$head.Text = ""

# ---------------------------------------------------------------------
# Slides 16-18 : recadrage des captures d'écran MCD + fond blanc
# rendu transparent (clrChange sur le blip)
# ---------------------------------------------------------------------
$s16 = $p.Slides.Item(16)
$shp = $s16.Shapes.Item("Picture 2")
$shp.Left   = 432000 / $EMU
$shp.Top    = 927056 / $EMU
$shp.Width  = 8542240 / $EMU
$shp.Height = 5580000 / $EMU
$shp.PictureFormat.TransparencyColor = 16777215

$s17 = $p.Slides.Item(17)
$shp = $s17.Shapes.Item("Picture 2")
$shp.Left = 432000 / $EMU
$shp.Top  = 928800 / $EMU
$shp.PictureFormat.TransparencyColor = 16777215

$s18 = $p.Slides.Item(18)
$shp = $s18.Shapes.Item("Picture 2")
$shp.Left   = 432480 / $EMU
$shp.Top    = 927054 / $EMU
$shp.Width  = 8542237 / $EMU
$shp.Height = 5580000 / $EMU
$shp.PictureFormat.TransparencyColor = 16777215
